$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-covering range: copy formatting for new rows 47-50 from row 46
$ws.Range("A46:B46").Copy()
$ws.Range("A47:B50").PasteSpecial(-4122)

# Write cell values for rows 13-50 (shifted / new primer rows)
$ws.Cells.Item(13, 1).Value = 'aceE_del-1'
$ws.Cells.Item(13, 2).Value = 'GCGTCACAGACATGAAATTGGT'
$ws.Cells.Item(14, 1).Value = 'aceE_del-9'
$ws.Cells.Item(14, 2).Value = 'AGCCATTATTCTTTTACCTCGGGTTATTCCTTATCTATCT'
$ws.Cells.Item(15, 1).Value = 'aceE_del-3'
$ws.Cells.Item(15, 2).Value = 'TTAACACCAAACTCGCGTGC'
$ws.Cells.Item(16, 1).Value = 'aceE_del-10'
$ws.Cells.Item(16, 2).Value = 'GAGGTAAAAGAATAATGGCTATCGA'
$ws.Cells.Item(17, 1).Value = 'aceE_del-2'
$ws.Cells.Item(17, 2).Value = 'TCGGCTTGAAAGGTTTGCACGGGTTATTCCTTATCTATCT'
$ws.Cells.Item(18, 1).Value = 'aceE_del-4'
$ws.Cells.Item(18, 2).Value = 'ACTACCTCATTACTGGATGCGAGGTAAAAGAATAATGGCTATCGA'
$ws.Cells.Item(19, 1).Value = 'aceE_del-7'
$ws.Cells.Item(19, 2).Value = 'GTGCAAACCTTTCAAGCCGA'
$ws.Cells.Item(20, 1).Value = 'aceE_del-8'
$ws.Cells.Item(20, 2).Value = 'GCATCCAGTAATGAGGTAGTTGC'
$ws.Cells.Item(21, 1).Value = 'test-aceE_del-1'
$ws.Cells.Item(21, 2).Value = 'CCGGAAACTCCGCTGGGCGA'
$ws.Cells.Item(22, 1).Value = 'test-aceE_del-2'
$ws.Cells.Item(22, 2).Value = 'GGAGCTGCTTCTGCACGTTT'
$ws.Cells.Item(23, 1).Value = 'test-aceE_del-3'
$ws.Cells.Item(23, 2).Value = 'ACGTAAAGTCTACATTTGTGCA'
$ws.Cells.Item(24, 1).Value = 'test-aceE_del-4'
$ws.Cells.Item(24, 2).Value = 'GGAGCTGCTTCTGCACGTTT'
$ws.Cells.Item(25, 1).Value = 'pntA_promoter_sub-1'
$ws.Cells.Item(25, 2).Value = 'CGAGGTTTGTGCCGTAAAGC'
$ws.Cells.Item(26, 1).Value = 'pntA_promoter_sub-9'
$ws.Cells.Item(26, 2).Value = 'AGGTATAATGCTAGCACGAATCTAGAGAAAGATTGGACGTACCATAATGCGAATTGGCATACCAAGAG'
$ws.Cells.Item(27, 1).Value = 'pntA_promoter_sub-3'
$ws.Cells.Item(27, 2).Value = 'ACTTGGTGATGCGGTAGTCG'
$ws.Cells.Item(28, 1).Value = 'pntA_promoter_sub-10'
$ws.Cells.Item(28, 2).Value = 'TTCGTGCTAGCATTATACCTAGGACTGAGCTAGCTGTCAAGGCGCGGTGATAGTGGGATAAACACCT'
$ws.Cells.Item(29, 1).Value = 'pntA_promoter_sub-2'
$ws.Cells.Item(29, 2).Value = 'TCGGCTTGAAAGGTTTGCACATGCGAATTGGCATACCAAGAG'
$ws.Cells.Item(30, 1).Value = 'pntA_promoter_sub-4'
$ws.Cells.Item(30, 2).Value = 'ACTACCTCATTACTGGATGCGTGATAGTGGGATAAACACCT'
$ws.Cells.Item(31, 1).Value = 'pntA_promoter_sub-7'
$ws.Cells.Item(31, 2).Value = 'GTGCAAACCTTTCAAGCCGA'
$ws.Cells.Item(32, 1).Value = 'pntA_promoter_sub-8'
$ws.Cells.Item(32, 2).Value = 'GCATCCAGTAATGAGGTAGTTGC'
$ws.Cells.Item(33, 1).Value = 'test-pntA_promoter_sub-1'
$ws.Cells.Item(33, 2).Value = 'CCGGAAACTCCGCTGGGCGA'
$ws.Cells.Item(34, 1).Value = 'test-pntA_promoter_sub-2'
$ws.Cells.Item(34, 2).Value = 'TGTCGAACGGGACCATCATC'
$ws.Cells.Item(35, 1).Value = 'test-pntA_promoter_sub-3'
$ws.Cells.Item(35, 2).Value = 'TAATTTCGCCCGCACGGAT'
$ws.Cells.Item(36, 1).Value = 'test-pntA_promoter_sub-4'
$ws.Cells.Item(36, 2).Value = 'TGTCGAACGGGACCATCATC'
$ws.Cells.Item(37, 1).Value = 'Cgl1452_ins-1'
$ws.Cells.Item(37, 2).Value = 'CACTGCGCGGGATTTTATGG'
$ws.Cells.Item(38, 1).Value = 'Cgl1452_ins-9'
$ws.Cells.Item(38, 2).Value = 'TCAATACTCTTTTTGGCGCGCATGTGAACGCCTGACCAGG'
$ws.Cells.Item(39, 1).Value = 'Cgl1452_ins-3'
$ws.Cells.Item(39, 2).Value = 'CGATGTCGCTGGCGTTAATG'
$ws.Cells.Item(40, 1).Value = 'Cgl1452_ins-10'
$ws.Cells.Item(40, 2).Value = 'GCTCCGAGGTTGAAGCTTAAGCATCCGGCATGAACAAAGC'
$ws.Cells.Item(41, 1).Value = 'Cgl1452_ins-2'
$ws.Cells.Item(41, 2).Value = 'TCGGCTTGAAAGGTTTGCACCATGTGAACGCCTGACCAGG'
$ws.Cells.Item(42, 1).Value = 'Cgl1452_ins-4'
$ws.Cells.Item(42, 2).Value = 'ACTACCTCATTACTGGATGCGCATCCGGCATGAACAAAGC'
$ws.Cells.Item(43, 1).Value = 'Cgl1452_ins-7'
$ws.Cells.Item(43, 2).Value = 'GTGCAAACCTTTCAAGCCGA'
$ws.Cells.Item(44, 1).Value = 'Cgl1452_ins-8'
$ws.Cells.Item(44, 2).Value = 'GCATCCAGTAATGAGGTAGTTGC'
$ws.Cells.Item(45, 1).Value = 'Cgl1452_ins-5'
$ws.Cells.Item(45, 2).Value = 'CGCGCCAAAAAGAGTATTGACT'
$ws.Cells.Item(46, 1).Value = 'Cgl1452_ins-6'
$ws.Cells.Item(46, 2).Value = 'TTAAGCTTCAACCTCGGAGCG'
$ws.Cells.Item(47, 1).Value = 'test-Cgl1452_ins-1'
$ws.Cells.Item(47, 2).Value = 'CCGGAAACTCCGCTGGGCGA'
$ws.Cells.Item(48, 1).Value = 'test-Cgl1452_ins-2'
$ws.Cells.Item(48, 2).Value = 'TGACTTGTTAGCCGGTCAGC'
$ws.Cells.Item(49, 1).Value = 'test-Cgl1452_ins-3'
$ws.Cells.Item(49, 2).Value = 'AGTCGCTAAAGTCAGGCCAT'
$ws.Cells.Item(50, 1).Value = 'test-Cgl1452_ins-4'
$ws.Cells.Item(50, 2).Value = 'TGACTTGTTAGCCGGTCAGC'
